$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.695.53'
$ws.Range("E2").Value = '  -0.60%  '

$ws.Range("D3").Value = '1.597.41'
$ws.Range("E3").Value = '  -1.11%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("E5").Value = '  -0.33%  '

$ws.Range("E6").Value = '  +0.21%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -0.76%  '

$ws.Range("E9").Value = '  -1.53%  '

$ws.Range("D10").Value = '''19.77'
$ws.Range("E10").Value = '  +0.00%  '

$ws.Range("D11").Value = '''0.0840'
$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("D12").Value = '1.820.86'
$ws.Range("E12").Value = '  -1.12%  '

$ws.Range("D13").Value = '1.595.06'
$ws.Range("E13").Value = '  -1.13%  '

$ws.Range("E14").Value = '  -1.29%  '

$ws.Range("E15").Value = '  -1.97%  '

$ws.Range("D16").Value = '''65.07'
$ws.Range("E16").Value = '  +1.79%  '

$ws.Range("D17").Value = '26.697.91'
$ws.Range("E17").Value = '  -0.52%  '

$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").Value = '''210.21'
$ws.Range("E19").Value = '  -0.16%  '

$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").Value = '''6.76'

$ws.Range("E22").Value = '  -0.78%  '

$ws.Range("E23").Value = '  -1.58%  '

$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").Value = '''146.71'
$ws.Range("E25").Value = '  +0.20%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  -4.14%  '

$ws.Range("E28").Value = '  +1.91%  '

$ws.Range("D29").Value = '''15.32'
$ws.Range("E29").Value = '  -0.51%  '

$ws.Range("E30").Value = '  -0.06%  '

$ws.Range("E31").Value = '  -0.49%  '

$ws.Range("E32").Value = '  -1.58%  '

$ws.Range("D33").Value = '''0.667'
$ws.Range("E33").Value = '  -6.02%  '

$ws.Range("E34").Value = '  -1.70%  '

$ws.Range("D35").Value = '1.296.87'
$ws.Range("E35").Value = '  -2.28%  '

$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("E37").Value = '  -4.16%  '

$ws.Range("D38").Value = '''0.0172'
$ws.Range("E38").Value = '  -1.41%  '

$ws.Range("D39").Value = '''0.844'
$ws.Range("E39").Value = '  +1.74%  '

$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''5.38'
$ws.Range("E41").Value = '  +1.47%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.790'
$ws.Range("E42").Value = '  -0.53%  '

$ws.Range("E43").Value = '  -0.54%  '

$ws.Range("D44").Value = '''63.78'
$ws.Range("E44").Value = '  +0.32%  '

$ws.Range("D45").Value = '1.733.81'
$ws.Range("E45").Value = '  -1.11%  '

$ws.Range("D46").Value = '''0.885'
$ws.Range("E46").Value = '  +7.31%  '

$ws.Range("D47").Value = '''90.15'
$ws.Range("E47").Value = '  +0.95%  '

$ws.Range("E48").Value = '  +0.55%  '

$ws.Range("D49").Value = '''0.100'
$ws.Range("E49").Value = '  +1.90%  '

$ws.Range("E50").Value = '  -1.46%  '

$ws.Range("D51").Value = '''7.51'
$ws.Range("E51").Value = '  +0.30%  '
